$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "우리 학교에 오지 말아야 할 분"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/not-wanted-students/#utm_source=rss&utm_medium=rss&utm_campaign=not-wanted-students"

$ws.Range("D51").Value = "[PySide6] DeprecationWarning: Fuction: 'globalPos() const' is marked as deprecated 경고 메시지 출력 안되게 하기"
$ws.Range("E51").Value = "https://bskyvision.com/1230"
